$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.416.52"
$ws.Range("E2").Value = "'  -3.91%  "
$ws.Range("D3").Value = "'3.295.00"
$ws.Range("E3").Value = "'  -6.25%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'549.56"
$ws.Range("E5").Value = "'  -4.33%  "
$ws.Range("D6").Value = "'171.96"
$ws.Range("E6").Value = "'  -3.81%  "
$ws.Range("D7").Value = "'0.608"
$ws.Range("E7").Value = "'  -4.24%  "
$ws.Range("E8").Value = "'  -0.10%  "
$ws.Range("D9").Value = "'3.290.14"
$ws.Range("E9").Value = "'  -6.23%  "
$ws.Range("D10").Value = "'0.616"
$ws.Range("E10").Value = "'  -2.70%  "
$ws.Range("E11").Value = "'  +0.28%  "
$ws.Range("D12").Value = "'53.46"
$ws.Range("E12").Value = "'  -2.61%  "
$ws.Range("D13").Value = "'0.0000269"
$ws.Range("E13").Value = "'  -1.58%  "
$ws.Range("D14").Value = "'8.95"
$ws.Range("E14").Value = "'  -2.95%  "
$ws.Range("D15").Value = "'3.816.73"
$ws.Range("E15").Value = "'  -6.41%  "
$ws.Range("D16").Value = "'18.09"
$ws.Range("E16").Value = "'  -1.52%  "
$ws.Range("E17").Value = "'  -3.70%  "
$ws.Range("D18").Value = "'3.290.97"
$ws.Range("E18").Value = "'  -6.29%  "
$ws.Range("D19").Value = "'11.69"
$ws.Range("E19").Value = "'  -3.70%  "
$ws.Range("D20").Value = "'63.230.98"
$ws.Range("E20").Value = "'  -4.08%  "
$ws.Range("D21").Value = "'0.963"
$ws.Range("E21").Value = "'  -4.35%  "
$ws.Range("D22").Value = "'422.51"
$ws.Range("E22").Value = "'  +1.74%  "
$ws.Range("D23").Value = "'4.60"
$ws.Range("E23").Value = "'  +5.01%  "
$ws.Range("D24").Value = "'4.06"
$ws.Range("E24").Value = "'  -2.79%  "
$ws.Range("D25").Value = "'83.79"
$ws.Range("E25").Value = "'  -2.16%  "
$ws.Range("D26").Value = "'13.08"
$ws.Range("E26").Value = "'  +1.97%  "
$ws.Range("D27").Value = "'10.55"
$ws.Range("E27").Value = "'  -3.49%  "
$ws.Range("D28").Value = "'2.80"
$ws.Range("E28").Value = "'  -1.86%  "
$ws.Range("D29").Value = "'8.63"
$ws.Range("E29").Value = "'  -4.20%  "
$ws.Range("D30").Value = "'29.35"
$ws.Range("E30").Value = "'  -3.35%  "
$ws.Range("D31").Value = "'6.57"
$ws.Range("E31").Value = "'  +2.35%  "
$ws.Range("D32").Value = "'587.45"
$ws.Range("E32").Value = "'  -5.69%  "
$ws.Range("D33").Value = "'11.33"
$ws.Range("E33").Value = "'  -2.80%  "
$ws.Range("E34").Value = "'  -3.77%  "
$ws.Range("D35").Value = "'58.01"
$ws.Range("E35").Value = "'  -2.72%  "
$ws.Range("E36").Value = "'  -0.19%  "
$ws.Range("E37").Value = "'  -6.79%  "
$ws.Range("B38").Value = "'InjectiveProtocol"
$ws.Range("C38").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'35.09"
$ws.Range("E38").Value = "'  -5.98%  "
$ws.Range("B39").Value = "'Stacks"
$ws.Range("C39").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'3.41"
$ws.Range("E39").Value = "'  +2.30%  "
$ws.Range("D40").Value = "'0.0₃0739"
$ws.Range("E40").Value = "'  -8.33%  "
$ws.Range("D41").Value = "'0.362"
$ws.Range("E41").Value = "'  -4.60%  "
$ws.Range("B42").Value = "'Maker"
$ws.Range("C42").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'3.076.86"
$ws.Range("E42").Value = "'  -6.04%  "
$ws.Range("B43").Value = "'FirstDigitalUSD"
$ws.Range("C43").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "'  -0.01%  "
$ws.Range("D44").Value = "'2.80"
$ws.Range("E44").Value = "'  -3.91%  "
$ws.Range("E45").Value = "'  -2.51%  "
$ws.Range("D46").Value = "'0.0402"
$ws.Range("E46").Value = "'  -4.09%  "
$ws.Range("E47").Value = "'  -3.47%  "
$ws.Range("D48").Value = "'0.129"
$ws.Range("E48").Value = "'  -3.04%  "
$ws.Range("D49").Value = "'2.57"
$ws.Range("E49").Value = "'  -5.62%  "
$ws.Range("D50").Value = "'132.74"
$ws.Range("E50").Value = "'  -4.39%  "
$ws.Range("E51").Value = "'  -5.13%  "
